# "taking latest code changes"
# On the "Test Cases" sheet, flip the Runmode column (C) from "N" to "Y"
# for the first five test rows so that they all run, and move the active
# cell selection to C6 to reflect where editing last left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("C2:C6").Value = "Y"

$ws.Range("C6").Select()
